# Applies the cryptos list price/volume refresh described by the commit
# "Updated cryptos list on Fri Apr  5 16:45:13 UTC 2024 with GitHub Actions".
# Only column D (Price) and column E (Volume(1h)) text values change; all other
# cells (rank, coin name, link) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.805.48"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.319.14"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'581.96"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "'174.35"
$ws.Range("E6").Value = "  -6.67%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.583"
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("D9").Value = "3.315.65"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("E10").Value = "  -4.92%  "
$ws.Range("D11").Value = "'0.575"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").Value = "'45.30"
$ws.Range("E12").Value = "  -4.86%  "
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").Value = "'665.92"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").Value = "3.856.93"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("D17").Value = "67.917.77"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "3.315.33"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("D20").Value = "'17.43"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").Value = "'10.88"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").Value = "'0.887"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").Value = "'5.37"
$ws.Range("E23").Value = "  +4.90%  "
$ws.Range("D24").Value = "'16.93"
$ws.Range("E24").Value = "  -5.89%  "
$ws.Range("D25").Value = "'97.57"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("E26").Value = "  -5.20%  "
$ws.Range("D27").Value = "'2.67"
$ws.Range("E27").Value = "  -6.74%  "
$ws.Range("D28").Value = "'9.25"
$ws.Range("E28").Value = "  -5.55%  "
$ws.Range("D29").Value = "'33.55"
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("D30").Value = "'8.39"
$ws.Range("E30").Value = "  -3.80%  "
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("D32").Value = "'583.59"
$ws.Range("E32").Value = "  -5.08%  "
$ws.Range("D33").Value = "'10.95"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("D35").Value = "3.720.83"
$ws.Range("E35").Value = "  -7.47%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "'57.54"
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").Value = "'3.31"
$ws.Range("E38").Value = "  -13.92%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'32.47"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("D41").Value = "'2.62"
$ws.Range("E41").Value = "  -7.02%  "
$ws.Range("D42").Value = "'3.08"
$ws.Range("E42").Value = "  -5.24%  "
$ws.Range("E43").Value = "  -3.61%  "
$ws.Range("D44").Value = "0.0₃0662"
$ws.Range("E44").Value = "  -6.14%  "
$ws.Range("D45").Value = "'3.27"
$ws.Range("E45").Value = "  -4.37%  "
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("D47").Value = "'2.58"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'1.34"
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("D51").Value = "'127.78"
$ws.Range("E51").Value = "  -0.41%  "
